$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Copy existing cell formatting onto the new H:J columns ---
# Header row (row 1) uses the same style as the existing header cells (e.g. B1).
$ws.Range("B1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# Data rows (rows 2-4) use the same style as the existing data cells (e.g. B2).
$ws.Range("B2").Copy()
$ws.Range("H2:J4").PasteSpecial(-4122)

# --- New header values ---
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# --- New data values: report date, legislator name, legislator id ---
# The date is written with a leading apostrophe so Excel keeps it as literal
# text ("2012-04-23") instead of auto-converting it to a date serial value.
$ws.Cells.Item(2, 8).Value = "'2012-04-23"
$ws.Cells.Item(2, 9).Value = "李鴻鈞"
$ws.Cells.Item(2, 10).Value = 898

$ws.Cells.Item(3, 8).Value = "'2012-04-23"
$ws.Cells.Item(3, 9).Value = "李鴻鈞"
$ws.Cells.Item(3, 10).Value = 898

$ws.Cells.Item(4, 8).Value = "'2012-04-23"
$ws.Cells.Item(4, 9).Value = "李鴻鈞"
$ws.Cells.Item(4, 10).Value = 898
